# Automatische test-sync: 2025-06-27 22:14:50
#
# Adds a new log row (row 4) for the "Openingstijden / Locatie" test mail to
# the "Logs" sheet, adds the matching aggregate row to the "Dashboard" sheet,
# extends the conditional formatting ranges that covered rows 2:3 to cover
# 2:4, and updates the bar chart's category/value series references so the
# chart picks up the new Dashboard row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet - append row 4
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(4, 1).Value = "Wanneer zijn jullie open?"
$logs.Cells.Item(4, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(4, 3).Value = "Testmail #1: Wanneer zijn jullie open?"
$logs.Cells.Item(4, 4).Value = "Openingstijden / Locatie"

$antwoord = @"
Beste klant,
Bedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.
Met vriendelijke groet,
[Naam bedrijf]
"@
$logs.Cells.Item(4, 5).Value = $antwoord

$logs.Cells.Item(4, 6).Value = "2025-06-27 22:14:36"
$logs.Cells.Item(4, 7).Value = "Ja"
$logs.Cells.Item(4, 8).Value = "Nee"
$logs.Cells.Item(4, 9).Value = "Ja"

# ---------------------------------------------------------------------
# 2. Logs sheet - extend conditional formatting (rows 2:3 -> 2:4)
# ---------------------------------------------------------------------
function Extend-FormatConditions($range, $newRange) {
    $fcs = $range.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

Extend-FormatConditions $logs.Range("D2:D3") $logs.Range("D2:D4")
Extend-FormatConditions $logs.Range("G2:G3") $logs.Range("G2:G4")
Extend-FormatConditions $logs.Range("H2:H3") $logs.Range("H2:H4")
Extend-FormatConditions $logs.Range("I2:I3") $logs.Range("I2:I4")

# ---------------------------------------------------------------------
# 3. Dashboard sheet - append aggregate row 4
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(4, 1).Value = "Openingstijden / Locatie"
$dash.Cells.Item(4, 2).Value = 1

# ---------------------------------------------------------------------
# 4. Update the chart's category/value series so it spans rows 2:4
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)

$series.XValues = "='Dashboard'!`$A`$2:`$A`$4"
$series.Values = "='Dashboard'!`$B`$2:`$B`$4"
